# Update cryptocurrency price values in column D
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = "249.35"
    3 = "22.97"
    4 = "5.418"
    5 = "0.05636"
    6 = "3.422"
    7 = "6.358"
    8 = "0.8127"
    9 = "0.9142"
    10 = "0.1439"
    11 = "0.07532"
    12 = "0.03098"
    13 = "0.03111"
    14 = "0.09337"
    15 = "3.558"
    16 = "0.001590"
    17 = "0.04757"
    18 = "0.0005789"
    19 = "0.006392"
    20 = "0.005002"
    26 = "0.1296"
    28 = "0.0003031"
    40 = "0.04041"
    41 = "0.006811"
    42 = "0.1071"
    43 = "0.002720"
    44 = "0.007555"
    45 = "0.00005791"
    47 = "0.4998"
    48 = "0.2410"
    49 = "0.00002099"
}

foreach ($row in $updates.Keys) {
    $cell = $ws.Range("D$row")
    # Force text storage so formatted numeric strings (leading/trailing zeros) are preserved exactly
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$row]
    $cell.Style = "Normal"
}
